# Weekly update: a new day's price record for "Ajo" (Femacal de La Calera)
# is inserted above the existing historical rows. Insert a whole new row at
# sheet row 156, which pushes the former rows 156..219 down to 157..220,
# and populate the newly-opened row 156 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data (old rows 156-219) down by one row.
$ws.Rows.Item(156).Insert()

# Fill in the new record in the now-empty row 156.
$ws.Cells.Item(156, 1).Value = 3
$ws.Cells.Item(156, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 44466
$ws.Cells.Item(156, 5).Value = 5
$ws.Cells.Item(156, 6).Value = 100112003
$ws.Cells.Item(156, 7).Value = "Ajo"
$ws.Cells.Item(156, 8).Value = "Chino"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 75
$ws.Cells.Item(156, 11).Value = 16000
$ws.Cells.Item(156, 12).Value = 17000
$ws.Cells.Item(156, 13).Value = 16467
$ws.Cells.Item(156, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(156, 15).Value = "China"
$ws.Cells.Item(156, 16).Value = 1647
$ws.Cells.Item(156, 17).Value = 10
$ws.Cells.Item(156, 18).Value = "Hortaliza"
